# capital_adequacy_2022_Q2.xlsx — add the "Table_2" coefficients sheet and
# drop the leftover empty placeholder cells on "Table_1".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Table_1: remove the stray empty inline-string cells (A3, B2, B6, ...)
#    These used to hold an empty string; the edit drops them entirely so
#    the row no longer carries a cell at that reference.
# ---------------------------------------------------------------------
$emptyCells = @(
    "B2", "A3", "B6", "A7", "B8", "B9", "A10", "B11", "B12", "A13", "B14",
    "B25", "A26", "B27", "B34", "B35", "B37", "B38", "B48", "A49", "B50",
    "B51", "B52"
)
foreach ($addr in $emptyCells) {
    $ws1.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------
# 2) Add the new "Table_2" worksheet right after "Table_1"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# ---------------------------------------------------------------------
# 3) Fill in the header row + three data rows. Numeric / percent-looking
#    text ("6.0%", "60000", ...) must stay literal text, not be coerced
#    to a number, so we briefly force a text number format for those
#    cells, write the value, then clear the format back off again so no
#    extra cell style lingers in the saved file.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$headers = @("Əmsal", "Norma (Sistem əhəmiyyətli)", "Norma (Banklar istisna)", "Fakt")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$rows = @(
    @("9.  I dərəcəli  kapitalın  adekvatlıq əmsalı", "6.0%", "5.0%", "9.2%"),
    @("10. məcmu kapitalın  adekvatlıq  əmsalı", "12.0%", "10.0%", "13.3%"),
    @("11. Leverec əmsalı", "minimum 5%", "minimum 4%", "5.4%")
)

$r = 2
foreach ($row in $rows) {
    Set-TextValue $ws2.Cells.Item($r, 1) $row[0]
    Set-TextValue $ws2.Cells.Item($r, 2) $row[1]
    Set-TextValue $ws2.Cells.Item($r, 3) $row[2]
    Set-TextValue $ws2.Cells.Item($r, 4) $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4) Reuse Table_1!A1's style (bold, bordered, centered) for the header
#    row so we don't introduce a brand-new cell style. Done last so the
#    value-writes above don't disturb the copied format.
# ---------------------------------------------------------------------
$ws1.Range("A1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Keep "Table_1" as the selected/active sheet, same as before the edit.
$ws1.Activate()
